# Automatische test-sync: 2025-06-27 22:36:50
# Append the new "Wanneer zijn jullie open?" test-mail row (row 8) to the
# "Logs" sheet, grow the conditional-formatting ranges to include it, and
# bump the "Openingstijden / Locatie" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- New row 8 on "Logs" ---------------------------------------------------
$logs.Range("A8").Value = "Wanneer zijn jullie open?"
$logs.Range("B8").Value = "mailmind.test@zohomail.eu"
$logs.Range("C8").Value = "Testmail #1: Wanneer zijn jullie open?"
$logs.Range("D8").Value = "Openingstijden / Locatie"
$logs.Range("E8").Value = "Beste klant,`nBedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Range("F8").Value = "2025-06-27 22:36:39"
$logs.Range("G8").Value = "Ja"
$logs.Range("H8").Value = "Nee"
$logs.Range("I8").Value = "Ja"

# --- Grow conditional-formatting ranges from row 7 to row 8 ----------------
$logs.Range("D2:D7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D8"))
$logs.Range("G2:G7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G8"))
$logs.Range("H2:H7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H8"))
$logs.Range("I2:I7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I8"))

# --- "Dashboard" tally: Openingstijden / Locatie count 4 -> 5 ---------------
$dash.Range("B2").Value = 5
